{"js": "const replacements = [\n  [\"59\u00d718=1062\", \"87\u00d720=1740\"],\n  [\"81\u00d747=3807\", \"54\u00d776=4104\"],\n  [\"37\u00d735=1295\", \"64\u00d769=4416\"],\n  [\"38\u00d782=3116\", \"82\u00d760=4920\"],\n  [\"38\u00d731=1178\", \"73\u00d716=1168\"],\n  [\"83\u00d744=3652\", \"85\u00d770=5950\"],\n  [\"87\u00d714=1218\", \"88\u00d792=8096\"],\n  [\"12\u00d722=264\", \"38\u00d759=2242\"],\n  [\"84\u00d716=1344\", \"29\u00d712=348\"],\n  [\"98\u00d796=9408\", \"82\u00d715=1230\"],\n  [\"71\u00d755=3905\", \"52\u00d744=2288\"],\n  [\"30\u00d729=870\", \"27\u00d729=783\"],\n  [\"65\u00d749=3185\", \"24\u00d763=1512\"],\n  [\"76\u00d727=2052\", \"86\u00d733=2838\"],\n  [\"55\u00d765=3575\", \"83\u00d788=7304\"],\n  [\"72\u00d785=6120\", \"56\u00d735=1960\"],\n  [\"15\u00d791=1365\", \"28\u00d786=2408\"],\n  [\"80\u00d798=7840\", \"87\u00d793=8091\"],\n  [\"48\u00d723=1104\", \"30\u00d793=2790\"],\n  [\"29\u00d728=812\", \"46\u00d753=2438\"],\n  [\"68\u00d799=6732\", \"33\u00d781=2673\"],\n  [\"30\u00d782=2460\", \"76\u00d799=7524\"],\n  [\"26\u00d754=1404\", \"66\u00d779=5214\"],\n  [\"47\u00d745=2115\", \"47\u00d716=752\"],\n  [\"81\u00d782=6642\", \"19\u00d756=1064\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old='59\u00d718=1062'; new='87\u00d720=1740'},\n    @{old='81\u00d747=3807'; new='54\u00d776=4104'},\n    @{old='37\u00d735=1295'; new='64\u00d769=4416'},\n    @{old='38\u00d782=3116'; new='82\u00d760=4920'},\n    @{old='38\u00d731=1178'; new='73\u00d716=1168'},\n    @{old='83\u00d744=3652'; new='85\u00d770=5950'},\n    @{old='87\u00d714=1218'; new='88\u00d792=8096'},\n    @{old='12\u00d722=264'; new='38\u00d759=2242'},\n    @{old='84\u00d716=1344'; new='29\u00d712=348'},\n    @{old='98\u00d796=9408'; new='82\u00d715=1230'},\n    @{old='71\u00d755=3905'; new='52\u00d744=2288'},\n    @{old='30\u00d729=870'; new='27\u00d729=783'},\n    @{old='65\u00d749=3185'; new='24\u00d763=1512'},\n    @{old='76\u00d727=2052'; new='86\u00d733=2838'},\n    @{old='55\u00d765=3575'; new='83\u00d788=7304'},\n    @{old='72\u00d785=6120'; new='56\u00d735=1960'},\n    @{old='15\u00d791=1365'; new='28\u00d786=2408'},\n    @{old='80\u00d798=7840'; new='87\u00d793=8091'},\n    @{old='48\u00d723=1104'; new='30\u00d793=2790'},\n    @{old='29\u00d728=812'; new='46\u00d753=2438'},\n    @{old='68\u00d799=6732'; new='33\u00d781=2673'},\n    @{old='30\u00d782=2460'; new='76\u00d799=7524'},\n    @{old='26\u00d754=1404'; new='66\u00d779=5214'},\n    @{old='47\u00d745=2115'; new='47\u00d716=752'},\n    @{old='81\u00d782=6642'; new='19\u00d756=1064'}\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($r.old, $false, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
